# Apply the crypto price/volume refresh described in the commit:
# "Updated cryptos list on Tue Oct 17 09:00:15 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.470.83'
$ws.Range("E2").Value = '  +3.44%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.588.52'
$ws.Range("E3").Value = '  +1.36%  '

# Row 4
$ws.Range("E4").Value = '  +1.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.61'
$ws.Range("E5").Value = '  +0.78%  '

# Row 6
$ws.Range("E6").Value = '  +0.80%  '

# Row 7
$ws.Range("E7").Value = '  +1.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.47'
$ws.Range("E8").Value = '  +8.27%  '

# Row 9
$ws.Range("E9").Value = '  +0.36%  '

# Row 10
$ws.Range("E10").Value = '  +0.76%  '

# Row 11
$ws.Range("E11").Value = '  +1.96%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.815.24'
$ws.Range("E12").Value = '  +1.35%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.584.11'
$ws.Range("E13").Value = '  +1.18%  '

# Row 14
$ws.Range("E14").Value = '  +2.01%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.74'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.460.09'
$ws.Range("E16").Value = '  +3.44%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.10'
$ws.Range("E17").Value = '  +1.13%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.51'
$ws.Range("E18").Value = '  +2.07%  '

# Row 19
$ws.Range("E19").Value = '  -0.22%  '

# Row 20
$ws.Range("E20").Value = '  +0.44%  '

# Row 21
$ws.Range("E21").Value = '  +1.04%  '

# Row 22
$ws.Range("E22").Value = '  -1.24%  '

# Row 23
$ws.Range("E23").Value = '  -0.70%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  +0.60%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.76'
$ws.Range("E25").Value = '  +1.21%  '

# Row 26
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.56'
$ws.Range("E27").Value = '  -0.78%  '

# Row 28
$ws.Range("E28").Value = '  -0.69%  '

# Row 29
$ws.Range("E29").Value = '  +1.05%  '

# Row 30
$ws.Range("E30").Value = '  -0.80%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0472'
$ws.Range("E31").Value = '  +0.29%  '

# Row 32
$ws.Range("E32").Value = '  +0.45%  '

# Row 33
$ws.Range("E33").Value = '  +0.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.399.34'
$ws.Range("E34").Value = '  -3.28%  '

# Row 35
$ws.Range("E35").Value = '  -0.84%  '

# Row 36
$ws.Range("E36").Value = '  -9.02%  '

# Row 37
$ws.Range("E37").Value = '  +1.09%  '

# Row 38
$ws.Range("E38").Value = '  -0.60%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.55'
$ws.Range("E39").Value = '  +7.57%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.542'
$ws.Range("E40").Value = '  +0.10%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("E41").Value = '  -0.33%  '

# Row 42
$ws.Range("E42").Value = '  +1.04%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.57'
$ws.Range("E43").Value = '  -2.73%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.85'
$ws.Range("E44").Value = '  +0.73%  '

# Row 46
$ws.Range("E46").Value = '  -0.72%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.725.39'
$ws.Range("E47").Value = '  +1.29%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.22'
$ws.Range("E48").Value = '  +0.64%  '

# Row 49
$ws.Range("B49").Value = 'mCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.13'
$ws.Range("E49").Value = '  +1.45%  '

# Row 50
$ws.Range("E50").Value = '  +12.80%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0518'
$ws.Range("E51").Value = '  -1.34%  '
